$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: change Assigned To (C2) and fill in Opportunity Purpose (E2) / Opportunity Status (F2)
$ws.Range("C2").Value = "John Smith"
$ws.Range("E2").Value = "Scholarship"
$ws.Range("F2").Value = "Stewardship"

# Row 3: change Assigned To (C3) and fill in Opportunity Status (F3)
$ws.Range("C3").Value = "John Smith"
$ws.Range("F3").Value = "Solicitation"

# Row 4: change Assigned To (C4)
$ws.Range("C4").Value = "John Smith"

# Row 5: change Assigned To (C5)
$ws.Range("C5").Value = "John Smith"
